# Progress update commit: "adding new progress as of date 04 nov 2025"
#
# For every data row (3-24) on the "Training Dashboard" sheet:
#   - column H ("PERIOD TO EXPIRE") drops by 1 day
#   - column I ("LAST UPDATE") moves from 03-Nov-2025 to 04-Nov-2025
#
# Column I holds the date as literal text (not a real Excel date), so we
# can't just assign the string straight onto a General-formatted cell --
# Excel's COM layer auto-coerces a date-shaped string into a date serial.
# Instead we build the literal text in a scratch cell (as a formula result,
# which Excel never re-parses as a date) and PasteSpecial only the VALUE
# into the target cell, leaving its existing number format/style untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Training Dashboard")

$oldDate = "03-Nov-2025"
$newDate = "04-Nov-2025"

$scratch = $ws.Cells.Item(100, 100)

for ($row = 3; $row -le 24; $row++) {
    $hCell = $ws.Cells.Item($row, 8)   # PERIOD TO EXPIRE
    $iCell = $ws.Cells.Item($row, 9)   # LAST UPDATE

    if ($iCell.Value2 -eq $oldDate) {
        $hCell.Value2 = $hCell.Value2 - 1

        $scratch.Formula = '="' + $newDate + '"'
        $scratch.Copy()
        $iCell.PasteSpecial(-4163)   # xlPasteValues - value only, keeps $iCell's style
    }
}

$scratch.Clear()
$excel.CutCopyMode = $false
